$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column V holds the "vendedor" (seller) field, added as an FK-like
# reference column to the clients sheet.

# Match the header look-and-feel used by the other "secondary" headers
# (F1 contato, I1 celular, M1 complemento, T1 observacoes) by copying
# their formatting onto the new header cell before writing its text.
$ws.Range("T1").Copy($ws.Range("V1"))

# Header
$ws.Range("V1").Value = "vendedor"

# Sample/help row (row 2 is the "field description" example row in this
# template) gets the explanation text; row 3 is left blank for this column.
$ws.Range("V2").Value = "Nome do vendedor igual ao do cadastro de vendedores"

# Widen the new column to fit the longer help text.
$ws.Columns("V").ColumnWidth = 47

# Move the view/selection the way the author left it after adding the column.
[void]$ws.Range("V7").Select()
$excel.ActiveWindow.ScrollColumn = 13
